$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode (column C) to "Y" and Results (column D) to "SKIP" for all
# test case rows (2-11) - "Running all the suites".
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = "Y"
    $ws.Cells.Item($r, 4).Value = "SKIP"
}

# Update the active selection to reflect the refreshed Runmode column.
$ws.Activate()
$ws.Range("C2:C11").Select()
